$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" '44.845.95'
Set-TextValue "E2" '  +1.93%  '
Set-TextValue "D3" '2.251.91'
Set-TextValue "E3" '  +0.61%  '
Set-TextValue "E4" '  +0.34%  '
Set-TextValue "D5" '307.71'
Set-TextValue "E5" '  +0.75%  '
Set-TextValue "D6" '95.95'
Set-TextValue "E6" '  +1.62%  '
Set-TextValue "D7" '0.574'
Set-TextValue "E7" '  +0.65%  '
Set-TextValue "E8" '  +0.19%  '
Set-TextValue "D9" '0.524'
Set-TextValue "E9" '  +0.35%  '
Set-TextValue "D10" '35.56'
Set-TextValue "E10" '  +2.36%  '
Set-TextValue "E11" '  -0.08%  '
Set-TextValue "D12" '7.27'
Set-TextValue "E12" '  +1.09%  '
Set-TextValue "E13" '  +0.09%  '
Set-TextValue "D14" '2.595.06'
Set-TextValue "E14" '  +0.60%  '
Set-TextValue "D15" '2.297.60'
Set-TextValue "E15" '  +0.13%  '
Set-TextValue "D16" '0.843'
Set-TextValue "E16" '  +2.69%  '
Set-TextValue "D17" '13.61'
Set-TextValue "E17" '  +0.98%  '
Set-TextValue "D18" '44.546.75'
Set-TextValue "E18" '  +1.52%  '
Set-TextValue "D19" '0.0₃0956'
Set-TextValue "E19" '  -0.63%  '
Set-TextValue "D20" '12.09'
Set-TextValue "E20" '  -0.59%  '
Set-TextValue "D21" '6.34'
Set-TextValue "E21" '  +0.62%  '
Set-TextValue "D22" '65.60'
Set-TextValue "E22" '  +0.13%  '
Set-TextValue "D23" '238.76'
Set-TextValue "E23" '  +0.99%  '
Set-TextValue "E24" '  +2.46%  '
Set-TextValue "E25" '  +2.09%  '
Set-TextValue "E26" '  -0.18%  '
Set-TextValue "D27" '2.28'
Set-TextValue "E27" '  +3.33%  '
Set-TextValue "D28" '9.91'
Set-TextValue "E28" '  +0.63%  '
Set-TextValue "D29" '37.42'
Set-TextValue "E29" '  -1.62%  '
Set-TextValue "D30" '6.04'
Set-TextValue "E30" '  +0.17%  '
Set-TextValue "D31" '20.02'
Set-TextValue "E31" '  +0.80%  '
Set-TextValue "D32" '152.70'
Set-TextValue "D33" '0.0802'
Set-TextValue "E33" '  +0.51%  '
Set-TextValue "E34" '  +1.66%  '
Set-TextValue "D35" '3.11'
Set-TextValue "E35" '  -1.78%  '
Set-TextValue "E36" '  +0.49%  '
Set-TextValue "E37" '  -0.46%  '
Set-TextValue "E38" '  +6.08%  '
Set-TextValue "D39" '14.96'
Set-TextValue "E39" '  -0.05%  '
Set-TextValue "D40" '3.42'
Set-TextValue "E40" '  +1.38%  '
Set-TextValue "D41" '3.86'
Set-TextValue "E41" '  +1.11%  '
Set-TextValue "D42" '0.0306'
Set-TextValue "E42" '  +3.60%  '
Set-TextValue "E43" '  +0.31%  '
Set-TextValue "D44" '1.837.65'
Set-TextValue "E44" '  +5.69%  '
Set-TextValue "E45" '  +18.11%  '
Set-TextValue "D47" '79.96'
Set-TextValue "E47" '  -5.64%  '
Set-TextValue "B48" 'Aave'
Set-TextValue "C48" 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue "D48" '99.67'
Set-TextValue "E48" '  -0.13%  '
Set-TextValue "B49" 'ordi'
Set-TextValue "C49" 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
Set-TextValue "D49" '70.73'
Set-TextValue "E49" '  +2.85%  '
Set-TextValue "E50" '  -0.25%  '
Set-TextValue "E51" '  +2.64%  '

Write-Host "Applied all changes"
